# Add 2022-Q3 data
# --------------------------------------------------------------------------
# 1) Insert a brand-new "2022-Q3" worksheet right after "总计" by copying the
#    existing "2022-Q2" sheet (this preserves sheetPr/outlinePr/pageSetUpPr
#    and column/row formatting exactly), then overwrite its data with the
#    fresh 2022-Q3 fund-holding figures and drop the now-unused 15th row.
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3 above
#    the existing rows and shift everything else down by one row.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$totalWs = $wb.Worksheets.Item("总计")
$q2Ws    = $wb.Worksheets.Item("2022-Q2")

# ---- Step 1: create the new "2022-Q3" sheet from a copy of "2022-Q2" -----
$q2Ws.Copy([System.Reflection.Missing]::Value, $totalWs)
$q3Ws = $wb.Worksheets.Item(2)
$q3Ws.Name = "2022-Q3"

# Make sure text-like numeric columns (B:G) stay text, like the rest of the
# workbook, before writing the new values into them.
$q3Ws.Range("B2:G14").NumberFormat = "@"

# Row-by-row fund holdings for 2022-Q3 (code, name, size, equity position,
# position weight, holding value, position rank)
$q3Data = @(
    @("010365", "鹏华港股通中证香港银行投资指数（LOF）C",       "4.92", "94.13", "4.42", "0.2175", 8),
    @("501025", "鹏华港股通中证香港银行投资指数（LOF）A",       "2.38", "94.13", "4.42", "0.1052", 8),
    @("006810", "泰康港股通中证香港银行投资指数C",             "0.73", "94.14", "4.46", "0.0326", 8),
    @("006809", "泰康港股通中证香港银行投资指数A",             "0.69", "94.14", "4.46", "0.0308", 8),
    @("501305", "汇添富中证港股通高股息投资指数（LOF）A",       "0.87", "92.21", "3.33", "0.0290", 10),
    @("513530", "华泰柏瑞中证港股通高股息投资ETF（QDII）",      "0.78", "95.80", "3.46", "0.0270", 10),
    @("501310", "华宝标普沪港深中国增强价值指数（LOF）A",       "0.89", "93.50", "3.00", "0.0267", 7),
    @("159726", "华夏恒生中国内地企业高股息率ETF",             "0.84", "96.48", "2.59", "0.0218", 10),
    @("517900", "招商中证银行AH价格优选ETF",                  "0.30", "96.90", "5.16", "0.0155", 6),
    @("501306", "汇添富中证港股通高股息投资指数（LOF）C",       "0.17", "92.21", "3.33", "0.0057", 10),
    @("501307", "银河中证沪港深高股息指数（LOF）A",            "0.15", "90.33", "1.34", "0.0020", 7),
    @("007397", "华宝标普沪港深中国增强价值指数（LOF）C",       "0.04", "93.50", "3.00", "0.0012", 7),
    @("501308", "银河中证沪港深高股息指数（LOF）C",            "0.01", "90.33", "1.34", "0.0001", 7)
)

$r = 2
foreach ($row in $q3Data) {
    $q3Ws.Range("A" + $r).Value = ($r - 2)
    $q3Ws.Range("B" + $r).Value = $row[0]
    $q3Ws.Range("C" + $r).Value = $row[1]
    $q3Ws.Range("D" + $r).Value = $row[2]
    $q3Ws.Range("E" + $r).Value = $row[3]
    $q3Ws.Range("F" + $r).Value = $row[4]
    $q3Ws.Range("G" + $r).Value = $row[5]
    $q3Ws.Range("H" + $r).Value = $row[6]
    $r = $r + 1
}

# The copied sheet had 14 data rows (dimension A1:H15); 2022-Q3 only needs
# 13, so remove the now-superfluous last row entirely.
$q3Ws.Rows(15).Delete()

# ---- Step 2: update the "总计" summary sheet ------------------------------
# Push the existing 7 quarters down one row and place the new 2022-Q3
# summary figures in row 2.
$totalWs.Rows(2).Insert()
$totalWs.Range("B2:D2").ClearFormats()
$totalWs.Range("A2").Value = 0
$totalWs.Range("A2").Style = "Normal"

$totalWs.Range("B2").Value = "2022-Q3"
$totalWs.Range("C2").Value = 13
$totalWs.Range("D2").Value = 0.52

# Re-number the row index column (A) sequentially (0-based) and re-apply the
# header-matching style used by every other cell in column A.
$quarterRows = @(2, 3, 4, 5, 6, 7, 8, 9)
for ($i = 0; $i -lt $quarterRows.Length; $i++) {
    $cell = $totalWs.Range("A" + $quarterRows[$i])
    $cell.Value = $i
}
